$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Remove-HyperlinkAt($addr) {
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
            break
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

# Row 14 used to carry hyperlinks on D14 and F14 - drop them before clearing the cells.
Remove-HyperlinkAt '$D$14'
Remove-HyperlinkAt '$F$14'

# Row 13: the K'IAM / SRMARTRH entry is removed - clear C13, D13, F13 (keep A13, B13, E13).
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("F13").ClearContents()

# Row 14: the SmartRH entry is removed - clear C14, D14, F14 (keep A14).
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("F14").ClearContents()

# Update the active selection to F13 (matches the new selection recorded in the sheet).
$ws.Range("F13").Select()
